$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "BgITlKb1"

$ws.Range("A6").Value = "Alex"
$ws.Range("B6").Value = "sobrenome"
$ws.Range("C6").Value = "192.168.100.10"
$ws.Range("D6").Value = "alex.junio@fgvjr.com"
$ws.Range("E6").Value = "Oll0LYyL"

$ws.Range("A7").Value = "B"
$ws.Range("B7").Value = "R"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "1123456543234567654345678"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "brunoluiszrosa@gmail.com"
$ws.Range("E7").Value = "17SMAQqc"
